$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 193.748574183173
$ws.Range("D2").Value = 0.000000000000000000000000000000000000000082927827769872
$ws.Range("B3").Value = 4177.66974871083
$ws.Range("B4").Value = 1968.22914910464
